# Update "想去人数" (F column) figures across the four worksheets to match
# a freshly regenerated data snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 159
$ws.Range("F5").Value  = 427
$ws.Range("F6").Value  = 806
$ws.Range("F7").Value  = 239
$ws.Range("F8").Value  = 1156
$ws.Range("F9").Value  = 332
$ws.Range("F11").Value = 870
$ws.Range("F12").Value = 670
$ws.Range("F18").Value = 2890
$ws.Range("F19").Value = 2611
$ws.Range("F20").Value = 524
$ws.Range("F21").Value = 27
$ws.Range("F26").Value = 5237
$ws.Range("F27").Value = 589
$ws.Range("F29").Value = 20
$ws.Range("F31").Value = 294
$ws.Range("F32").Value = 1075

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 1109
$ws.Range("F9").Value  = 327
$ws.Range("F10").Value = 23
$ws.Range("F14").Value = 602
$ws.Range("F19").Value = 40
$ws.Range("F26").Value = 3895
$ws.Range("F30").Value = 197
$ws.Range("F31").Value = 47
$ws.Range("F35").Value = 9

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value  = 2434
$ws.Range("F6").Value  = 1027
$ws.Range("F9").Value  = 1300
$ws.Range("F10").Value = 351
$ws.Range("F11").Value = 95

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 2434
$ws.Range("F6").Value  = 1027
$ws.Range("F7").Value  = 1300
$ws.Range("F8").Value  = 351
$ws.Range("F9").Value  = 95
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 427
$ws.Range("F12").Value = 806
$ws.Range("F13").Value = 239
$ws.Range("F14").Value = 1156
$ws.Range("F15").Value = 332
$ws.Range("F16").Value = 870
$ws.Range("F17").Value = 670
$ws.Range("F18").Value = 1109
$ws.Range("F19").Value = 1109
$ws.Range("F24").Value = 2890
$ws.Range("F25").Value = 2611
$ws.Range("F26").Value = 524
$ws.Range("F27").Value = 27
$ws.Range("F30").Value = 5237
$ws.Range("F31").Value = 589
$ws.Range("F33").Value = 602
$ws.Range("F34").Value = 602
$ws.Range("F35").Value = 20
$ws.Range("F38").Value = 294
$ws.Range("F40").Value = 40
$ws.Range("F46").Value = 1075
$ws.Range("F47").Value = 197
$ws.Range("F48").Value = 47
